# Updated dissertation following all review comments
#
# 1) Bump the "datetimeFigureOut" date placeholder shown on every slide
#    master/layout from 4/15/2016 to 5/23/2016.
# 2) Rename the "Execution Time of Operation" label (inside the
#    "Amplitude:" box on slide 2) to "Response Time of Operation".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = "5/23/2016"
        }
    }
}

# Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Slide 2: "Execution Time of Operation" -> "Response Time of Operation"
# inside the "Rectangle 25" callout box (keeps the "Amplitude:" line as-is).
$slide = $p.Slides.Item(2)
$shapes = $slide.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $needle = "Execution Time of Operation"
        $idx = $full.IndexOf($needle)
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, $full.Length - $idx)
            $sub.Text = "Response Time of Operation"
        }
    }
}
